# FrontendTimeline.xlsx update:
#  - add a new timeline row (row 8): Day 5, Date 05/08/2025, 5 hours,
#    "Finished Mostly User Basic profile update, password reset, email change"
#  - total hours formula (D10) recalculates automatically
#  - leave the active selection on D18

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the row above (row 7) onto the new row 8 cells so the
# new entries keep the same styles (centered number / date format / etc.)
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)

$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)

$ws.Range("C7").Copy()
$ws.Range("C8").PasteSpecial(-4122)

$ws.Range("D7").Copy()
$ws.Range("D8").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Fill in the new row's data
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = 45874
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = "Finished Mostly User Basic profile update, password reset, email change"

# Match the saved selection from the source file
$ws.Range("D18").Select() | Out-Null
